# Update "provincias_spain" COVID-19 Spain worksheet:
#  - refresh the "Datos actualizados..." timestamp cell
#  - update case counts for the provinces/cities whose numbers changed
#    (the refreshed data set is still sorted by "Casos totales" descending,
#    which reshuffles a few rows' city labels along with their numbers)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 14:20"

# row, city, CasosTotales, CasosActivos, Recuperados, Muertes
$rows = @(
    @(19, "Malaga",              1349,   83, 1196, 70),
    @(20, "Gran Canaria",        1262,   57,  320, 11),
    @(21, "La Palma",            1262,   57,   47,  2),
    @(22, "Lanzarote",           1262,   57,   42,  3),
    @(23, "Fuerteventura",       1262,   57,   31,  0),
    @(24, "La Gomera",           1262,   57,    7,  0),
    @(25, "El Hierro",           1262,   57,    3,  0),
    @(29, "Sevilla",             1119,   17, 1066, 36),
    @(30, "Gipuzkoa/Guipuzcoa",  1113, 1796,  673, 44),
    @(31, "Caceres",             1067,   29,  924, 114),
    @(32, "Granada",             1061,   15,  979, 67),
    @(33, "Valladolid",           988,  187,  724, 77),
    @(34, "Leon",                 964,  176,  688, 100),
    @(39, "Cordoba",              661,    4,  642, 15),
    @(40, "Jaen",                 661,   17,  618, 26),
    @(41, "Guadalajara",          643,  296,  532, 97),
    @(42, "Castello/Castellon",   613,    9,  570, 34),
    @(46, "Cadiz",                539,   16,  509, 14),
    @(51, "Almeria",              251,    6,  229, 16),
    @(52, "Huesca",               244,   23,  207, 14),
    @(53, "Teruel",               236,   16,  205, 15),
    @(56, "Huelva",               177,    2,  171,  4)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
